# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets.
# Both sheets contain the same 28 data rows (rows 2-28), and the same set of
# rows in column F need to be bumped up to their new values.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 11726
    3  = 11354
    6  = 1027
    11 = 10770
    12 = 4163
    16 = 2468
    18 = 50
    19 = 131
    20 = 450
    21 = 11143
    22 = 10928
    28 = 14
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
